# Add 2022-Q3 data
# 1) Insert a new worksheet "2022-Q3" right after "总计" with the fund detail
#    data for that quarter.
# 2) Insert a new top row in "总计" summarizing the 2022-Q3 quarter, pushing
#    the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# header row + index column styling (bold, thin border, centered/top) to
# match the rest of the workbook's quarterly sheets
$q3HeaderStyle = $q3.Range("B1:H1")
$q3HeaderStyle.Font.Bold = $true
$q3HeaderStyle.Borders.LineStyle = 1
$q3HeaderStyle.HorizontalAlignment = -4108
$q3HeaderStyle.VerticalAlignment = -4160

$q3IndexStyle = $q3.Range("A2:A5")
$q3IndexStyle.Font.Bold = $true
$q3IndexStyle.Borders.LineStyle = 1
$q3IndexStyle.HorizontalAlignment = -4108
$q3IndexStyle.VerticalAlignment = -4160

$q3Rows = @(
    @(0, "519198", "万家颐和灵活配置混合A",  "9.55",  "93.93", "5.14", "0.4909", 10),
    @(1, "008979", "万家民丰回报一年持有期混合", "18.51", "29.35", "1.59", "0.2943", 9),
    @(2, "519197", "万家颐达灵活配置混合",    "2.25",  "45.36", "1.86", "0.0418", 9),
    @(3, "016620", "万家颐和灵活配置混合C",   "0.18",  "93.93", "5.14", "0.0093", 10)
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r,1).Value = $row[0]
    $q3.Cells.Item($r,2).Value = "'" + $row[1]
    $q3.Cells.Item($r,3).Value = $row[2]
    $q3.Cells.Item($r,4).Value = "'" + $row[3]
    $q3.Cells.Item($r,5).Value = "'" + $row[4]
    $q3.Cells.Item($r,6).Value = "'" + $row[5]
    $q3.Cells.Item($r,7).Value = "'" + $row[6]
    $q3.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Step 2: insert a new summary row for 2022-Q3 at the top of the "总计" sheet
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()

$summaryRows = @(
    @(0, "2022-Q3", 4,  0.84),
    @(1, "2022-Q2", 5,  1.06),
    @(2, "2022-Q1", 17, 1.86),
    @(3, "2021-Q4", 5,  1),
    @(4, "2021-Q3", 5,  0.67)
)

# Row insertion in this engine copies the format of the row above (the bold
# header) across the whole new row, which is correct for column A (bold
# index, matches every other row) but wrong for B:D (plain data cells) -
# clear that inherited formatting before filling B:D in.
$total.Range("B2:D2").ClearFormats()

$r = 2
foreach ($row in $summaryRows) {
    $total.Cells.Item($r,1).Value = $row[0]
    $total.Cells.Item($r,2).Value = $row[1]
    $total.Cells.Item($r,3).Value = $row[2]
    $total.Cells.Item($r,4).Value = $row[3]
    $r = $r + 1
}

$totalIndexStyle = $total.Range("A2")
$totalIndexStyle.Font.Bold = $true
$totalIndexStyle.Borders.LineStyle = 1
$totalIndexStyle.HorizontalAlignment = -4108
$totalIndexStyle.VerticalAlignment = -4160
